$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1143.5714
$ws.Range("J19").Value = 1075.75
$ws.Range("L19").Value = 1075.75
$ws.Range("N19").Value = -1425.75
$ws.Range("H76").Value = 4371.5557
$ws.Range("I76").Value = 4098.8335
$ws.Range("J76").Value = 4917
$ws.Range("K76").Value = 4098.8335
$ws.Range("L76").Value = 4917
$ws.Range("M76").Value = -3783.8335
$ws.Range("N76").Value = -5547
$ws.Range("H79").Value = 4371.5557
$ws.Range("I79").Value = 4098.8335
$ws.Range("J79").Value = 4917
$ws.Range("K79").Value = 4098.8335
$ws.Range("L79").Value = 4917
$ws.Range("M79").Value = -3006.8335
$ws.Range("N79").Value = -7101
$ws.Range("H98").Value = 31253520
$ws.Range("I98").Value = 34486200
$ws.Range("J98").Value = 4286.3335
$ws.Range("K98").Value = 34486200
$ws.Range("L98").Value = 4286.3335
$ws.Range("M98").Value = -34484702
$ws.Range("N98").Value = -7282.3335
$ws.Range("H103").Value = 744.45
$ws.Range("I103").Value = 385.5
$ws.Range("J103").Value = 1103.4
$ws.Range("K103").Value = 1156.5
$ws.Range("L103").Value = 3310.2
$ws.Range("M103").Value = -570.5
$ws.Range("N103").Value = -4482.200000000001
$ws.Range("H112").Value = 6627.2812
$ws.Range("I112").Value = 1000
$ws.Range("J112").Value = 6808.8066
$ws.Range("K112").Value = 3000
$ws.Range("L112").Value = 20426.4198
$ws.Range("M112").Value = -1892
$ws.Range("N112").Value = -22642.4198
$ws.Range("H122").Value = 31253520
$ws.Range("I122").Value = 34486200
$ws.Range("J122").Value = 4286.3335
$ws.Range("K122").Value = 103458600
$ws.Range("L122").Value = 12859.0005
$ws.Range("M122").Value = -103456150
$ws.Range("N122").Value = -17759.0005
$ws.Range("H132").Value = 2221.8125
$ws.Range("I132").Value = 2221.8125
$ws.Range("K132").Value = 6665.4375
$ws.Range("M132").Value = -4135.4375
$ws.Range("H135").Value = 385150.44
$ws.Range("I135").Value = 385150.44
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 3466353.96
$ws.Range("L135").Value = 0
$ws.Range("M135").Value = -3463818.96
$ws.Range("N135").Value = $null
$ws.Range("H140").Value = 59998
$ws.Range("J140").Value = 59998
$ws.Range("L140").Value = 59998
$ws.Range("N140").Value = -70358
$ws.Range("H141").Value = 1937.3529
$ws.Range("I141").Value = 1937.3529
$ws.Range("K141").Value = 5812.0587
$ws.Range("M141").Value = -632.0587000000005

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3381798.8
$ws.Range("I32").Value = 3679530.8
$ws.Range("K32").Value = 3679530.8
$ws.Range("M32").Value = -3679243.8
$ws.Range("H110").Value = 1309.8889
$ws.Range("I110").Value = 1223.625
$ws.Range("K110").Value = 1223.625
$ws.Range("M110").Value = 821.375
$ws.Range("H122").Value = 10428.875
$ws.Range("I122").Value = 12867.883
$ws.Range("K122").Value = 38603.649
$ws.Range("M122").Value = -36153.649
$ws.Range("H138").Value = 79424
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").Value = $null

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 6097.981
$ws.Range("I134").Value = 3713.9714
$ws.Range("K134").Value = 11141.9142
$ws.Range("M134").Value = -8606.914199999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 5447.2856
$ws.Range("J16").Value = 7543.364
$ws.Range("L16").Value = 7543.364
$ws.Range("N16").Value = -8117.364
$ws.Range("H31").Value = 8904.761
$ws.Range("I31").Value = 3736.75
$ws.Range("J31").Value = 11661.033
$ws.Range("K31").Value = 3736.75
$ws.Range("L31").Value = 11661.033
$ws.Range("M31").Value = -3441.75
$ws.Range("N31").Value = -12251.033
$ws.Range("H34").Value = 8904.761
$ws.Range("I34").Value = 3736.75
$ws.Range("J34").Value = 11661.033
$ws.Range("K34").Value = 3736.75
$ws.Range("L34").Value = 11661.033
$ws.Range("M34").Value = -3534.75
$ws.Range("N34").Value = -12065.033
$ws.Range("H113").Value = 5447.2856
$ws.Range("J113").Value = 7543.364
$ws.Range("L113").Value = 7543.364
$ws.Range("N113").Value = -11883.364
$ws.Range("H132").Value = 4387.2324
$ws.Range("I132").Value = 1676.7
$ws.Range("K132").Value = 5030.1
$ws.Range("M132").Value = -2500.1

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H128").Value = 216665.33
$ws.Range("I128").Value = 216665.33
$ws.Range("K128").Value = 649995.99
$ws.Range("M128").Value = -645015.99
$ws.Range("H138").Value = 129411.125
$ws.Range("I138").Value = 147184.14
$ws.Range("K138").Value = 441552.42
$ws.Range("M138").Value = -436412.42
$ws.Range("H139").Value = 55194.5
$ws.Range("I139").Value = 93853.73
$ws.Range("K139").Value = 281561.19
$ws.Range("M139").Value = -276421.19
$ws.Range("H141").Value = 2369.2856
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").Value = $null

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H24").Value = 9990
$ws.Range("I24").Value = 0
$ws.Range("J24").Value = 9990
$ws.Range("K24").Value = 0
$ws.Range("L24").Value = 9990
$ws.Range("M24").Value = $null
$ws.Range("N24").Value = -10336
$ws.Range("H122").Value = 42699.816
$ws.Range("I122").Value = 57439.74
$ws.Range("J122").Value = 7692.5
$ws.Range("K122").Value = 172319.22
$ws.Range("L122").Value = 23077.5
$ws.Range("M122").Value = -169869.22
$ws.Range("N122").Value = -27977.5
$ws.Range("H132").Value = 3032.9707
$ws.Range("I132").Value = 3032.9707
$ws.Range("K132").Value = 9098.9121
$ws.Range("M132").Value = -6568.9121

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5983.107
$ws.Range("I7").Value = 4070.8572
$ws.Range("J7").Value = 7895.357
$ws.Range("K7").Value = 4070.8572
$ws.Range("L7").Value = 7895.357
$ws.Range("M7").Value = -3958.8572
$ws.Range("N7").Value = -8119.357
$ws.Range("H61").Value = 5264.7
$ws.Range("I61").Value = 3772.182
$ws.Range("K61").Value = 3772.182
$ws.Range("M61").Value = -3570.182
$ws.Range("H113").Value = 5264.7
$ws.Range("I113").Value = 3772.182
$ws.Range("K113").Value = 3772.182
$ws.Range("M113").Value = -1602.182
$ws.Range("H126").Value = 5983.107
$ws.Range("I126").Value = 4070.8572
$ws.Range("J126").Value = 7895.357
$ws.Range("K126").Value = 12212.5716
$ws.Range("L126").Value = 23686.071
$ws.Range("M126").Value = -9742.571599999999
$ws.Range("N126").Value = -28626.071

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H43").Value = 7500
$ws.Range("I43").Value = 7500
$ws.Range("J43").Value = 0
$ws.Range("K43").Value = 7500
$ws.Range("L43").Value = 0
$ws.Range("M43").Value = -7351
$ws.Range("N43").Value = $null
$ws.Range("H81").Value = 9528731
$ws.Range("I81").Value = 1853.3334
$ws.Range("K81").Value = 3706.6668
$ws.Range("M81").Value = -2645.6668
$ws.Range("H84").Value = 9528731
$ws.Range("I84").Value = 1853.3334
$ws.Range("K84").Value = 18533.334
$ws.Range("M84").Value = -13229.334
$ws.Range("H100").Value = 811.58826
$ws.Range("J100").Value = 1075.4445
$ws.Range("L100").Value = 2150.889
$ws.Range("N100").Value = -3232.889
$ws.Range("H132").Value = 16140141
$ws.Range("I132").Value = 20002464
$ws.Range("J132").Value = 47132
$ws.Range("K132").Value = 60007392
$ws.Range("L132").Value = 141396
$ws.Range("M132").Value = -60004862
$ws.Range("N132").Value = -146456
$ws.Range("H136").Value = 43524780
$ws.Range("I136").Value = 100000690
$ws.Range("J136").Value = 81770.92
$ws.Range("K136").Value = 300002070
$ws.Range("L136").Value = 245312.76
$ws.Range("M136").Value = -299999520
$ws.Range("N136").Value = -250412.76
